# LOT2040.docx edit script
# Applies the content changes described in the commit diff via Word COM Find/Replace.
# [char]11 (vertical tab) is used as the "manual line break" token that Word's
# Find/Replace engine maps to <w:br/> when building runs.

$d = $word.ActiveDocument
$vt = [char]11

function Replace-Text($old, $new) {
    $result = $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        throw "Find/Replace failed for: $old"
    }
}

# 1. Heading3 title
Replace-Text "Genetic Engineering" "Theoretical and Practical Genetic Engineering"

# 2. Créditos-aula
Replace-Text "Créditos-aula: 3" "Créditos-aula: 5"

# 3. Carga horária
Replace-Text "Carga horária: 45 h" "Carga horária: 75 h"

# 4. Ativação
Replace-Text "Ativação: 01/01/2025" "Ativação: 15/07/2025"

# 5. Objetivos (PT)
Replace-Text "Introduzir os princípios básicos da genética com apresentação minuciosa dos importantes conceitos da genética clássica, genética molecular e genômica associada às técnicas e aplicações da genética molecular." "Aulas teóricas e práticas sobre os fundamentos de Engenharia Genética aplicados à Biotecnologia."

# 6. Objetivos (EN, italic)
Replace-Text "Introduce fundamental genetic principles while thoroughly covering essential concepts in classical genetics, molecular genetics, and genomics, alongside the techniques and applications of molecular genetics." "Theoretical and practical lessons on the fundamentals of Genetic Engineering applied to Biotechnology."

# 7. Docente(s) list bullet -> objectives/lesson plan paragraph
$old7 = "8711290 - Elisson Antônio da Costa Romanel" + $vt + "8853480 - Tatiane da Franca Silva"
$new7 = "Capacitar o aluno a compreender os principais conceitos e técnicas envolvidos na manipulação genética de organismos, bem como desenvolver habilidades em técnicas de genética molecular com enfoque em aplicações biotecnológicas." + $vt + "Aulas teóricas: 1) Introdução à Genética, 2) Nucleotídeos e Estrutura do DNA, 3) Genes e Cromossomos, 4) Replicação do DNA, 5) Transcrição e Processamento do RNA, 6) Código genético e Tradução, 7) Vetores e clonagem de DNA, 8) Genômica, 9) Regulação da Expressão Gênica, 10) Elementos Genéticos Transponíveis, 11) Mutação, Reparo de DNA e Recombinação. " + $vt + $vt + "Aulas práticas: 1) Técnicas de extração de ácido nucléico, 2) Elaboração de gel de agarose, 3) Métodos de quantificação de ácido nucléico, 4) Uso de enzimas de restrição, 5) Corrida de eletroforese, 6) Desenho de iniciadores para PCR, 7) Amplificação de gene por PCR convencional, 8) Clonagem em vetor plasmidial, 9) Triagem de clones positivos, 10) Visita (viagem didática complementar) a uma empresa está prevista, conforme disponibilidade."
Replace-Text $old7 $new7

# 8. Programa resumido (PT) - note: the original text uses a Greek question mark
#    look-alike character (U+037E) instead of a plain semicolon after "genética"
$gq = [char]0x037E
$old8 = "Introdução à genética" + $gq + " Estrutura e Replicação Molecular do DNA; Transcrição, Tradução e Código Genético; Mutação e Reparo do DNA; Regulação da Expressão Gênica; Genômica e Bioinformática; Técnicas e Aplicações da Genética Molecular."
$new8 = 'Notas - N distribuído no semestre. A composição das "N" fica a critério do docente. O curso consistirá em aulas expositivas, sessões de laboratório, discussões e exercícios. Projetos e/ou atividades que envolvam a preparação e/ou apresentações de seminários também estão previstos.'
Replace-Text $old8 $new8

# 9. Programa resumido (EN, italic)
Replace-Text "Introduction to genetics; Structure and Molecular Replication of DNA; Transcription, Translation and Genetic Code; DNA Mutation and Repair; Regulation of Gene Expression; Genomics and Bioinformatics; Techniques and Applications of Molecular Genetics." "Provide students with the knowledge to understand key concepts and techniques in the genetic manipulation of organisms, while developing practical skills in molecular genetics techniques with a focus on biotechnological."

# 10. Programa (PT)
Replace-Text "1. Introdução à Genética. 2. Nucleotídeos e Estrutura do DNA 3. Genes e Cromossomos 4. Replicação do DNA 5. Transcrição e Processamento do RNA 6. Código genético e Tradução 7. Clonagem de DNA 8. Genômica 9. Bioinformática 10. Regulação da Expressão Gênica 11. Elementos Genéticos Transponíveis 12. Mutação, Reparo de DNA e Recombinação 13. Técnicas e Aplicações da Genética Molecular." "MF = (somatório de N)/número de N (adequando o valor de N, quando houver peso distinto para as Ns)."

# 11. Programa (EN, italic) -> split theoretical/practical classes
$old11 = "1. Introduction to Genetics. 2. Nucleotides and DNA Structure 3. Genes and Chromosomes 4. DNA Replication 5. Transcription and RNA Processing 6. Genetic Code and Translation 7. DNA Cloning 8. Genomics 9. Bioinformatics 10. Regulation of Gene Expression 11. Elements Transposable Genetics 12. Mutation, DNA Repair and Recombination 13. Techniques and Applications of Molecular Genetics"
$new11 = "Theoretical classes: 1) Introduction to Genetics, 2) Nucleotides and DNA Structure, 3) Genes and Chromosomes, 4) DNA Replication, 5) Transcription and RNA Processing, 6) Genetic Code and Translation, 7) Vectors and DNA Cloning, 8) Genomics, 9) Gene Expression Regulation, 10) Transposable Genetic Elements, 11) Mutation, DNA Repair, and Recombination." + $vt + $vt + "Practical classes: 1) Techniques for nucleic acid extraction, 2) Preparation of agarose gels, 3) Nucleic acid quantification methods, 4) Use of restriction enzymes, 5) Gel electrophoresis, 6) Primer design for PCR, 7) Gene amplification using conventional PCR, 8) Cloning into plasmid vectors, 9) Screening for positive clones, and 10). A visit (complementary educational trip) to a company is planned, subject to availability."
Replace-Text $old11 $new11

# 12. Avaliação / Método text
Replace-Text 'Notas - N distribuído no semestre. A composição das "N" fica critério do docente.' "NF = (MF + PR)/2, onde PR é uma prova de recuperação. Prova de Recuperação (PR) para alunos com Média Final (MF) maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final (NF) igual ou maior do que 5,0."

# 13. Avaliação / Critério text -> full bibliography
Replace-Text "MF = MF = (somatório de N)/número de N (adequando o valor de N, quando houver peso distinto para as Ns)" "-Fundamentos de Genética. Peter Snustad e Michael Simmons (2017, 7ª Edição). Editora Guanabara Koogan. -Principles of Genetics. Dr. Peter Snustad and Michael Simmons (2019, 7th Edition). John Wiley and Sons. -Genética: Um enfoque conceitual. Benjamin Pierce (2016, 5ª Edição). Editora Guanabara Koogan. – Práticas e protocolos básicos de Biologia Molecular. Fernanda Matias (2021, 1ª Edição). Editora Blucher. -Genetics: A Conceptual Approach. Benjamin Pierce (2019, 7th Edition). W. H. Freeman. -Introdução à genética. Griffiths, Doebley, Peichel e Wassarman (2022 – 12ª Edição). Guanabara Koogan. -An Introduction to Genetic Analysis. Anthony Griffiths, John Doebley, Catherine Peichel, David A. Wassarman (2020 12th Edition). W. H. Freeman. -Molecular Biotechnology, Principles and Applications of Recombinant DNA. Bernard R. Glick and Cheryl L. Patten (2022, Sixth Edition). ASP Press."

# 14. Avaliação / Norma de recuperação text -> docente 1
Replace-Text "Norma de Recuperação: NF = (MF + PR)/2, onde PR é uma prova de recuperação. Prova de Recuperação (PR) para alunos com Média Final (MF) maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final (NF) igual ou maior do que 5,0." "8711290 - Elisson Antônio da Costa Romanel"

# 15. Old Bibliografia paragraph (3 runs separated by manual line breaks) -> docente 2
$old15 = "-Fundamentos de Genética. Peter Snustad e Michael Simmons (2017, 7ª Edição). Editora Guanabara Koogan. -Principles of Genetics. Dr. Peter Snustad and Michael Simmons (2019, 7th Edition). John Wiley and Sons. " + $vt + "-Genética: Um enfoque conceitual. Benjamin Pierce (2016, 5ª Edição). Editora Guanabara Koogan. -Genetics: A Conceptual Approach. Benjamin Pierce (2019, 7th Edition). W. H. Freeman. " + $vt + "-Introdução à genética. Griffiths, Doebley, Peichel e Wassarman (2022 – 12ª Edição). Guanabara Koogan. -An Introduction to Genetic Analysis. Anthony Griffiths, John Doebley, Catherine Peichel, David A. Wassarman (2020 12th Edition). W. H. Freeman."
Replace-Text $old15 "8853480 - Tatiane da Franca Silva"

"Done"
